# Auto-generated Excel COM-interop edit script
# Applies cached-value updates to the Leve profit calculation columns (H:N)
# across the ALC, ARM, CRP, GSM, and WVR sheets, per the scheduled market-price refresh.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3946.6667
$ws.Cells.Item(113, 9).Value = 3762.5
$ws.Cells.Item(113, 10).Value = 4157.143
$ws.Cells.Item(113, 11).Value = 3762.5
$ws.Cells.Item(113, 12).Value = 4157.143
$ws.Cells.Item(113, 13).Value = -508.5
$ws.Cells.Item(113, 14).Value = -10665.143
$ws.Cells.Item(132, 8).Value = 5954818
$ws.Cells.Item(132, 9).Value = 7521118
$ws.Cells.Item(132, 10).Value = 2879.8
$ws.Cells.Item(132, 11).Value = 22563354
$ws.Cells.Item(132, 12).Value = 8639.400000000001
$ws.Cells.Item(132, 13).Value = -22560824
$ws.Cells.Item(132, 14).Value = -13699.4

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1182.762
$ws.Cells.Item(45, 9).Value = 1222.4
$ws.Cells.Item(45, 11).Value = 1222.4
$ws.Cells.Item(45, 13).Value = -845.4000000000001
$ws.Cells.Item(74, 8).Value = 970.30304
$ws.Cells.Item(74, 9).Value = 942.1539
$ws.Cells.Item(74, 10).Value = 1074.8572
$ws.Cells.Item(74, 11).Value = 942.1539
$ws.Cells.Item(74, 12).Value = 1074.8572
$ws.Cells.Item(74, 13).Value = -68.15390000000002
$ws.Cells.Item(74, 14).Value = -2822.8572
$ws.Cells.Item(77, 8).Value = 970.30304
$ws.Cells.Item(77, 9).Value = 942.1539
$ws.Cells.Item(77, 10).Value = 1074.8572
$ws.Cells.Item(77, 11).Value = 4710.7695
$ws.Cells.Item(77, 12).Value = 5374.286
$ws.Cells.Item(77, 13).Value = -342.7695000000003
$ws.Cells.Item(77, 14).Value = -14110.286
$ws.Cells.Item(132, 8).Value = 4527.4287
$ws.Cells.Item(132, 9).Value = 4552.143
$ws.Cells.Item(132, 10).Value = 4403.857
$ws.Cells.Item(132, 11).Value = 13656.429
$ws.Cells.Item(132, 12).Value = 13211.571
$ws.Cells.Item(132, 13).Value = -11126.429
$ws.Cells.Item(132, 14).Value = -18271.571

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(129, 8).Value = 49749.5
$ws.Cells.Item(129, 9).Value = 49001
$ws.Cells.Item(129, 10).Value = 49999
$ws.Cells.Item(129, 11).Value = 49001
$ws.Cells.Item(129, 12).Value = 49999
$ws.Cells.Item(129, 13).Value = -44001
$ws.Cells.Item(129, 14).Value = -59999
$ws.Cells.Item(130, 8).Value = 78990
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 78990
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 78990
$ws.Cells.Item(130, 14).Value = -89030
$ws.Cells.Item(131, 8).Value = 30966.666
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 30966.666
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 30966.666
$ws.Cells.Item(131, 14).Value = -41046.666
$ws.Cells.Item(132, 8).Value = 3755.3684
$ws.Cells.Item(132, 9).Value = 3629.7144
$ws.Cells.Item(132, 10).Value = 3828.6667
$ws.Cells.Item(132, 11).Value = 10889.1432
$ws.Cells.Item(132, 12).Value = 11486.0001
$ws.Cells.Item(132, 13).Value = -8359.143199999999
$ws.Cells.Item(132, 14).Value = -16546.0001
$ws.Cells.Item(133, 8).Value = 25589.428
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 25589.428
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 25589.428
$ws.Cells.Item(133, 14).Value = -30649.428
$ws.Cells.Item(134, 8).Value = 859.2083
$ws.Cells.Item(134, 9).Value = 784.8605
$ws.Cells.Item(134, 10).Value = 1498.6
$ws.Cells.Item(134, 11).Value = 2354.5815
$ws.Cells.Item(134, 12).Value = 4495.799999999999
$ws.Cells.Item(134, 13).Value = 180.4184999999998
$ws.Cells.Item(134, 14).Value = -9565.8
$ws.Cells.Item(135, 8).Value = 42500
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 42500
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 42500
$ws.Cells.Item(135, 14).Value = -52640
$ws.Cells.Item(137, 8).Value = 29000
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 29000
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 29000
$ws.Cells.Item(137, 14).Value = -39200
$ws.Cells.Item(138, 8).Value = 37800
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 37800
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 37800
$ws.Cells.Item(138, 14).Value = -48080
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(140, 8).Value = 35000
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 35000
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 35000
$ws.Cells.Item(140, 14).Value = -45360
$ws.Cells.Item(141, 8).Value = 57500
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 57500
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 57500
$ws.Cells.Item(141, 14).Value = -67860

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(125, 8).Value = 40000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 40000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 40000
$ws.Cells.Item(125, 14).Value = -44920
$ws.Cells.Item(126, 8).Value = 1290.1428
$ws.Cells.Item(126, 9).Value = 1240.2222
$ws.Cells.Item(126, 10).Value = 1380
$ws.Cells.Item(126, 11).Value = 3720.6666
$ws.Cells.Item(126, 12).Value = 4140
$ws.Cells.Item(126, 13).Value = -1250.6666
$ws.Cells.Item(126, 14).Value = -9080
$ws.Cells.Item(127, 8).Value = 23500
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 23500
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 23500
$ws.Cells.Item(127, 14).Value = -33420
$ws.Cells.Item(128, 8).Value = 104799
$ws.Cells.Item(128, 9).Value = 250000
$ws.Cells.Item(128, 10).Value = 68498.75
$ws.Cells.Item(128, 11).Value = 250000
$ws.Cells.Item(128, 12).Value = 68498.75
$ws.Cells.Item(128, 13).Value = -245020
$ws.Cells.Item(128, 14).Value = -78458.75
$ws.Cells.Item(129, 8).Value = 41499.75
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 41499.75
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 41499.75
$ws.Cells.Item(129, 14).Value = -51499.75
$ws.Cells.Item(130, 8).Value = 58750
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 58750
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 58750
$ws.Cells.Item(130, 14).Value = -68790
$ws.Cells.Item(131, 8).Value = 27000.666
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 27000.666
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 27000.666
$ws.Cells.Item(131, 14).Value = -37080.666
$ws.Cells.Item(132, 8).Value = 46394.523
$ws.Cells.Item(132, 9).Value = 60886.59
$ws.Cells.Item(132, 10).Value = 5333.6665
$ws.Cells.Item(132, 11).Value = 182659.77
$ws.Cells.Item(132, 12).Value = 16000.9995
$ws.Cells.Item(132, 13).Value = -180129.77
$ws.Cells.Item(132, 14).Value = -21060.9995
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 8).Value = 23950
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 23950
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 71850
$ws.Cells.Item(134, 14).Value = -76920
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(136, 8).Value = 17345
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 17345
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 52035
$ws.Cells.Item(136, 14).Value = -57135
$ws.Cells.Item(137, 8).Value = 78000
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 78000
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 78000
$ws.Cells.Item(137, 14).Value = -88200
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(139, 8).Value = 25913
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 25913
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 25913
$ws.Cells.Item(139, 14).Value = -36193
$ws.Cells.Item(140, 8).Value = 44500
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 44500
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 44500
$ws.Cells.Item(140, 14).Value = -54860
$ws.Cells.Item(141, 8).Value = 37402.9
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 37402.9
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 37402.9
$ws.Cells.Item(141, 14).Value = -47762.9

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1557
$ws.Cells.Item(81, 9).Value = 1350.8572
$ws.Cells.Item(81, 10).Value = 3000
$ws.Cells.Item(81, 11).Value = 2701.7144
$ws.Cells.Item(81, 12).Value = 6000
$ws.Cells.Item(81, 13).Value = -1640.7144
$ws.Cells.Item(81, 14).Value = -8122
$ws.Cells.Item(84, 8).Value = 1557
$ws.Cells.Item(84, 9).Value = 1350.8572
$ws.Cells.Item(84, 10).Value = 3000
$ws.Cells.Item(84, 11).Value = 13508.572
$ws.Cells.Item(84, 12).Value = 30000
$ws.Cells.Item(84, 13).Value = -8204.572
$ws.Cells.Item(84, 14).Value = -40608

Write-Output "Applied all Leve profit updates."